# Add "actual value" (F) and "execution result" (G) columns for the white-box
# test / telegram-fee test rows, mirroring the expected value into the actual
# value column and marking every case as "pass".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$expected = @{
    2  = "650.0"
    3  = "655.0"
    4  = "990.0"
    5  = "995.0"
    6  = "586.0"
    7  = "592.0"
    8  = "1054.0"
    9  = "424.0"
    10 = "433.0"
    11 = "1216.0"
    12 = "1225.0"
    13 = "820.0"
    14 = "1060.0"
    15 = "0.0"
}

foreach ($row in 2..15) {
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $expected[$row]
    $fCell.ClearFormats()

    $gCell = $ws.Cells.Item($row, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = "pass"
    $gCell.ClearFormats()
}

$ws.Range("I6:I7").Select()
